# Auto-generated edit script for Economic Dashboard update (2025-11-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style swaps (toggle the yellow "freshly updated date" highlight) ---
# Donor cells that keep style 47 (no highlight) / 48 (yellow highlight) throughout this script
$donor47 = $ws.Range("C3")
$donor48 = $ws.Range("N3")

$donor48.Copy()
$ws.Range("C17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C17").Value = 45901

$donor48.Copy()
$ws.Range("C18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C18").Value = 45901

$donor47.Copy()
$ws.Range("C28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C28").Value = 45870

$donor47.Copy()
$ws.Range("C29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C29").Value = 45870

$donor47.Copy()
$ws.Range("C30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C30").Value = 45870

$donor47.Copy()
$ws.Range("C31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C31").Value = 45870

$donor48.Copy()
$ws.Range("C32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C32").Value = 45870

$donor48.Copy()
$ws.Range("C33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C33").Value = 45870

$donor48.Copy()
$ws.Range("C34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C34").Value = 45870

$donor48.Copy()
$ws.Range("N37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N37").Value = 45901

$donor48.Copy()
$ws.Range("N38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N38").Value = 45901

$donor48.Copy()
$ws.Range("N39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N39").Value = 45982

$excel.CutCopyMode = $false

# --- Value-only updates ---
$ws.Range("F17").Value = 0.001637846129888487
$ws.Range("G17").Value = 0.005952781743703284
$ws.Range("H17").Value = 0.006492096487988874
$ws.Range("I17").Value = 0.009680198742914703
$ws.Range("J17").Value = -0.007880419346928291
$ws.Range("F18").Value = 0.04261856857461769
$ws.Range("G18").Value = 0.05024109124453582
$ws.Range("H18").Value = 0.04134309243240536
$ws.Range("I18").Value = 0.0441678737351054
$ws.Range("J18").Value = 0.03367187567662758
$ws.Range("N29").Value = 45985
$ws.Range("R29").Value = 2.16
$ws.Range("S29").Value = 2.14
$ws.Range("U29").Value = 2.18
$ws.Range("N30").Value = 45985
$ws.Range("Q30").Value = 2.23
$ws.Range("S30").Value = 2.24
$ws.Range("U30").Value = 2.27
$ws.Range("F32").Value = -0.0007577910391209919
$ws.Range("G32").Value = -0.002144439521488684
$ws.Range("H32").Value = 0.00424479189840965
$ws.Range("I32").Value = -0.001038119019209582
$ws.Range("J32").Value = 0.001983365069813559
$ws.Range("F33").Value = 0.008358981150223696
$ws.Range("G33").Value = 0.01371833355505385
$ws.Range("H33").Value = 0.006650544135429186
$ws.Range("I33").Value = 0.00270366735076292
$ws.Range("J33").Value = 0.009949782230052044
$ws.Range("F34").Value = 75.84010000000001
$ws.Range("G34").Value = 75.9897
$ws.Range("H34").Value = 76.246
$ws.Range("I34").Value = 76.01730000000001
$ws.Range("J34").Value = 76.18980000000001
$ws.Range("Q37").Value = -0.002680233580537372
$ws.Range("R37").Value = -0.003492679271734134
$ws.Range("S37").Value = -0.002089419937226888
$ws.Range("T37").Value = 0.0006667672348770193
$ws.Range("U37").Value = 0.004670985029205044
$ws.Range("Q38").Value = 0.01291790739138332
$ws.Range("R38").Value = 0.01449043557431158
$ws.Range("S38").Value = 0.01630494859795866
$ws.Range("T38").Value = 0.01943770804710052
$ws.Range("U38").Value = 0.02364466296472451
$ws.Range("Q39").Value = 122.235
$ws.Range("R39").Value = 122.0735
$ws.Range("S39").Value = 121.8845
$ws.Range("T39").Value = 121.6042
$ws.Range("U39").Value = 121.5131
$ws.Range("N47").Value = 45982
$ws.Range("N48").Value = 45982
$ws.Range("Q48").Value = 3.51
$ws.Range("R48").Value = 3.55
$ws.Range("T48").Value = 3.58
$ws.Range("U48").Value = 3.6
$ws.Range("N49").Value = 45982
$ws.Range("Q49").Value = 3.62
$ws.Range("R49").Value = 3.68
$ws.Range("S49").Value = 3.71
$ws.Range("T49").Value = 3.7
$ws.Range("U49").Value = 3.72
$ws.Range("N50").Value = 45982
$ws.Range("Q50").Value = 4.06
$ws.Range("R50").Value = 4.1
$ws.Range("S50").Value = 4.13
$ws.Range("T50").Value = 4.12
$ws.Range("U50").Value = 4.13
$ws.Range("N52").Value = 45982
$ws.Range("Q52").Value = 5.88
$ws.Range("R52").Value = 5.9
$ws.Range("S52").Value = 5.92
$ws.Range("T52").Value = 5.91
$ws.Range("U52").Value = 5.9

# --- Restore unrelated empty numeric cells (G43:J43) that the host runtime's
# --- save round-trip otherwise coerces from blank to 0; these are untouched
# --- by this commit's diff, so keep them blank to match the target state.
$ws.Range("G43:J43").ClearContents()
